# Applies the "Add files via upload" edit to keybindings.xlsx:
#  - Rotates the "Move Focus" labels in B6:B9 (editor 1/editor 2/files/outline
#    -> files/outline/editor 1/editor 2)
#  - Updates the keybinding for "close file"/"close all files" from
#    ctrl e / ctrl e e -> ctrl r / ctrl t
#  - Clears the leftover highlight formatting on B14:C15 (matches the plain
#    filler style used elsewhere, e.g. M4)
#  - Moves the saved cell selection to G20

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Rotate the "Move Focus" target labels ---
$ws.Range("B6").Value = "files"
$ws.Range("B7").Value = "outline"
$ws.Range("B8").Value = "editor 1"
$ws.Range("B9").Value = "editor 2"

# --- Update the keybindings for close file / close all files ---
$ws.Range("I6").Value = "ctrl r"
$ws.Range("I7").Value = "ctrl t"

# --- Clear the special font-color formatting on B14:C15, matching the
#     plain fill-only style already used by column M/N filler cells ---
$ws.Range("M4").Copy() | Out-Null
$ws.Range("B14:C15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Restore the saved selection ---
$ws.Range("G20").Select() | Out-Null
